$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds the "last changed" date for every data row.
# The sheet's data occupies rows 2..556 (row 1 is the header).
# Bump the date from 45181 (2023-09-12) to 45182 (2023-09-13) for every
# data row, matching the recorded diff.
for ($r = 2; $r -le 556; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}
